# Add files via upload
# Appends 4 new rows (391-394) of variable/label_fr/label_en/color data to
# the "liste référence" sheet, matching the new "Jan 2nd / June 23rd"
# residential & total entries (EN + FR) introduced in the source workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A (variable) & Column C (label_en) -----------------------
# A and C carry the same English text throughout this sheet.
$ws.Range("A391").Value = "Jan 2nd, residential"
$ws.Range("A392").Value = "June 23rd, residential"
$ws.Range("A393").Value = "Jan 2nd, total"
$ws.Range("A394").Value = "June 23rd, total"

$ws.Range("C391").Value = "Jan 2nd, residential"
$ws.Range("C392").Value = "June 23rd, residential"
$ws.Range("C393").Value = "Jan 2nd, total"
$ws.Range("C394").Value = "June 23rd, total"

# --- Column B (label_fr) ----------------------------------------------
$ws.Range("B391").Value = "2 janv., résidentiel"
$ws.Range("B392").Value = "23 juin, résidentiel"
$ws.Range("B393").Value = "2 janv., total"
$ws.Range("B394").Value = "23 juin, total"

# --- Column D (color) ---------------------------------------------------
$ws.Range("D391").Value = "#9bdb9a"
$ws.Range("D392").Value = "#d20a11"
$ws.Range("D393").Value = "#5487a4"
$ws.Range("D394").Value = "#f6b4a4"

# Give the two new "residential" rows (391-392) a thin box border around
# their A/C variable + label_en cells, to set them apart as a new group,
# just like in the source edit. (Multi-area Range property sets only
# reliably hit the first area in this host, so apply per area.)
foreach ($addr in @("A391:A392", "C391:C392")) {
    $borderCells = $ws.Range($addr)
    $borderCells.HorizontalAlignment = -4131
    $borderCells.VerticalAlignment = -4160
    $borderCells.Font.Name = "Calibri"
    $borderCells.Font.Size = 11
    $borderCells.Borders.LineStyle = 1
}

# Copy the existing per-row "color swatch" formatting for column D from
# same-styled rows elsewhere in the sheet, so the new rows reuse the
# workbook's existing style entries instead of creating new ones.
$ws.Range("D275").Copy()
$ws.Range("D391").PasteSpecial(-4122)

$ws.Range("D57").Copy()
$ws.Range("D392").PasteSpecial(-4122)

$ws.Range("D2").Copy()
$ws.Range("D393").PasteSpecial(-4122)

$ws.Range("D17").Copy()
$ws.Range("D394").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Scroll / select near the new rows, mirroring where the author ended up.
$ws.Range("D394").Select()

Write-Output "done"
